# Add a new data row (row 34) to Sheet1 recording the analysis run for
# sg_rr_100_028, and update the sheet's view/selection state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row -----------------------------------------------------
# Columns (per header row 24):
#   A Data CSV filename              G prominence/dBm
#   B Wavelength step size/nm        H distance
#   C Start array index              I approx_fsr/nm
#   D End array index                J fsr_mean/nm
#   E Start wavelength/nm            K fsr_std error/nm
#   F End wavelength/nm              L double count check passed?
#                                     M (notes on prominence selection)

$ws.Range("A34").Value = "sg_rr_100_028 2023-12-08 16-58-05.csv"
$ws.Range("B34").Value = 0.01
$ws.Range("C34").Value = 1000
$ws.Range("D34").Value = 5001
$ws.Range("E34").Value = 1530
$ws.Range("F34").Value = 1570
$ws.Range("G34").Value = 0.5
$ws.Range("H34").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I34").Value = 1.7
$ws.Range("J34").Value = 0.98274999999999801
$ws.Range("K34").Value = [double]"3.2814064370514399E-3"
# Note column (M) is entered before the check column (L) so the new
# shared-string entries land in the same order as the source workbook.
$ws.Range("M34").Value = "prominence left at 0.5, as height span of noisy parts of data looks roughly around this value roughly at its max span looking by eye at a glance."
$ws.Range("L34").Value = "yes"

# --- View state ---------------------------------------------------------
# Scroll so row 22 is at the top of the visible area, and move the active
# selection down to the cell below the freshly added row.
$win = $excel.ActiveWindow
$ws.Range("L35").Select()
$win.ScrollRow = 22
$win.ScrollColumn = 1

# Resize/position the workbook window to match the recorded state.
$win.Left = 720
$win.Top = 720
$win.Width = 14400
$win.Height = 9350

$wb.Save()
